$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'youth basketball leggings with knee pads'
$ws.Cells.Item(2, 1).Value = 'youth basketball compression pants with knee pads'
$ws.Cells.Item(3, 1).Value = 'softball sliding knee pads'
$ws.Cells.Item(4, 1).Value = 'baseball sliding knee pads'
$ws.Cells.Item(5, 1).Value = 'hex pads basketball knee youth'
$ws.Cells.Item(6, 1).Value = 'youth basketball knee pads'
$ws.Cells.Item(7, 1).Value = 'sliding knee pads softball'
$ws.Cells.Item(8, 1).Value = 'youth black knee pads for basketball'
$ws.Cells.Item(9, 1).Value = 'crossfit leggings'
$ws.Cells.Item(10, 1).Value = 'knee pads for basketball'
$ws.Cells.Item(11, 1).Value = 'basketball pads'
$ws.Cells.Item(12, 1).Value = 'basketball knee'
$ws.Cells.Item(13, 1).Value = 'basketball leggings for boys'
$ws.Cells.Item(14, 1).Value = 'basketball leg pads tights'
$ws.Cells.Item(15, 1).Value = 'compression basketball pants men'
$ws.Cells.Item(16, 1).Value = 'crossfit knee sleeves'
$ws.Cells.Item(17, 1).Value = 'padded leg sleeves for basketball youth'
$ws.Cells.Item(18, 1).Value = 'hex knee pad'
$ws.Cells.Item(19, 1).Value = 'compression pants men capri'
$ws.Cells.Item(20, 1).Value = 'knee sleeves for crossfit'
$ws.Cells.Item(21, 1).Value = 'basketball leg compression pants'
$ws.Cells.Item(22, 1).Value = 'compression leggings for men'
$ws.Cells.Item(23, 1).Value = 'snowboarding knee pads men'
$ws.Cells.Item(24, 1).Value = 'mens tights and leggings'
$ws.Cells.Item(25, 1).Value = 'adult football pants with pads'
$ws.Cells.Item(26, 1).Value = 'volleyball knee pads xxl'
$ws.Cells.Item(27, 1).Value = 'compression basketball pants boys'
$ws.Cells.Item(28, 1).Value = 'basketball knee compression sleeve'
$ws.Cells.Item(29, 1).Value = 'padded knee compression sleeve'
$ws.Cells.Item(30, 1).Value = 'knee pads basketball girls'
$ws.Cells.Item(31, 1).Value = 'knee pads for work pants'
$ws.Cells.Item(32, 1).Value = 'knee brace wrestling'
$ws.Cells.Item(33, 1).Value = 'crossfit knee support'
$ws.Cells.Item(34, 1).Value = 'men compression leggings'
$ws.Cells.Item(35, 1).Value = 'basketball training gear'
$ws.Cells.Item(36, 1).Value = 'sports knee pad'
$ws.Cells.Item(37, 1).Value = 'hex padded compression leg sleeve'
$ws.Cells.Item(38, 1).Value = 'knee sleeve weightlifting men'
$ws.Cells.Item(39, 1).Value = 'hex pads'
$ws.Cells.Item(40, 1).Value = 'cycling capris padded'
$ws.Cells.Item(41, 1).Value = 'knee pads volleyball men'
$ws.Cells.Item(42, 1).Value = 'boy compression pants basketball'
$ws.Cells.Item(43, 1).Value = 'crossfit knee compression sleeve'
$ws.Cells.Item(44, 1).Value = 'crossfit knee sleeve'
$ws.Cells.Item(45, 1).Value = 'mens tights'
$ws.Cells.Item(46, 1).Value = 'protector de rodillas basketball'
$ws.Cells.Item(47, 1).Value = 'workout pads'
$ws.Cells.Item(48, 1).Value = 'mens athletic tights'
$ws.Cells.Item(49, 1).Value = 'athletic capri'
$ws.Cells.Item(50, 1).Value = 'mens basketball leggings'
$ws.Cells.Item(51, 1).Value = 'mens sliding pants baseball'
$ws.Cells.Item(52, 1).Value = 'padded compression pants'
$ws.Cells.Item(53, 1).Value = 'sliding pants baseball'
$ws.Cells.Item(54, 1).Value = 'tights capri men'
$ws.Cells.Item(55, 1).Value = 'tights with knee pads basketball'
$ws.Cells.Item(56, 1).Value = 'workout knee pads for men'
$ws.Cells.Item(57, 1).Value = 'workout sliding pads'
$ws.Cells.Item(58, 1).Value = 'wrestling knee pad sleeve'
$ws.Cells.Item(59, 1).Value = 'youth knee pads basketball'
$ws.Cells.Item(60, 1).Value = 'youth wrestling knee pads'
$ws.Cells.Item(61, 1).Value = 'leggings knee pads'
$ws.Cells.Item(62, 1).Value = 'xl knee pads'
$ws.Cells.Item(63, 1).Value = 'gel knee pads wrestling'
$ws.Cells.Item(64, 1).Value = 'men s leggings compression'
$ws.Cells.Item(65, 1).Value = 'capri tights for men'
$ws.Cells.Item(66, 1).Value = 'padded knee pads'
$ws.Cells.Item(67, 1).Value = 'football knee pads youth'
$ws.Cells.Item(68, 1).Value = 'knee compression sleeve cycling'
$ws.Cells.Item(69, 1).Value = 'volleyball knee pads for girls'
$ws.Cells.Item(70, 1).Value = 'youth football leggings'
$ws.Cells.Item(71, 1).Value = 'youth basketball leggings boys'
$ws.Cells.Item(72, 1).Value = 'wrestling knee sleeve'
$ws.Cells.Item(73, 1).Value = 'compression sleeve youth'
$ws.Cells.Item(74, 1).Value = 'workout capri pants'
$ws.Cells.Item(75, 1).Value = 'workout tights'
$ws.Cells.Item(76, 1).Value = 'basketball pants'
$ws.Cells.Item(77, 1).Value = 'compression leggings men'
$ws.Cells.Item(78, 1).Value = 'youth knee pads'
$ws.Cells.Item(79, 1).Value = 'basketball leggins with knee pads'
$ws.Cells.Item(80, 1).Value = 'knee pads black mens'
$ws.Cells.Item(81, 1).Value = 'basketball knee pads white'
$ws.Cells.Item(82, 1).Value = 'knee pads white basketball'
$ws.Cells.Item(83, 1).Value = 'black workout capris'
$ws.Cells.Item(84, 1).Value = 'black workout leggings capri'
$ws.Cells.Item(85, 1).Value = 'workout pads for knees'
$ws.Cells.Item(86, 1).Value = 'knee pads for basketball blue'
$ws.Cells.Item(87, 1).Value = 'knee tights for men'
$ws.Cells.Item(88, 1).Value = 'mcdavid padded leg sleeve'
$ws.Cells.Item(89, 1).Value = 'teen leggings'
$ws.Cells.Item(90, 1).Value = 'workout pants for men'
$ws.Cells.Item(91, 1).Value = 'mens training pants'
$ws.Cells.Item(92, 1).Value = 'basketball knee pads women'
$ws.Cells.Item(93, 1).Value = 'legging for men workout'
$ws.Cells.Item(94, 1).Value = 'mens athletic leggings'
$ws.Cells.Item(95, 1).Value = 'mens capri'
$ws.Cells.Item(96, 1).Value = 'workout support'
$ws.Cells.Item(97, 1).Value = 'athletic legging'
$ws.Cells.Item(98, 1).Value = 'basketball knee pads boys youth'
$ws.Cells.Item(99, 1).Value = 'basketball spandex knee pads'
$ws.Cells.Item(100, 1).Value = 'training gear six pad'
